$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pdfTarget = "http://www.hepforge.org/archive/lhapdf/pdfsets/6.1/NNPDF30_nlo_nf_5_pdfas.tar.gz"

# Rows 32-35: add the DY amcatnlo PDF-set / xml-ID columns (E-I), matching the
# pattern already used by the surrounding rows (e.g. row 28).
foreach ($r in 32..35) {
    $ws.Cells.Item($r, 5).Value = 292200                      # E: N events-ish id
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $pdfTarget)     # F: PDF set (+ hyperlink)
    $ws.Cells.Item($r, 6).Value = "NNPDF30_nlo_nf_5_pdfas"
    $ws.Cells.Item($r, 7).Value = "1001-1009"                 # G: ME/PS xml IDs
    $ws.Cells.Item($r, 8).Value = "2001-2100"                 # H: PDF Variation xml IDs
    $ws.Cells.Item($r, 9).Value = "2101-2102"                 # I: Strong Coupling xml IDs
}

# Reflect the view-state change captured in the saved workbook (user had
# scrolled/selected a different range while making the edit).
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Application.ActiveWindow.ScrollRow = 21
$ws.Range("G34:I35").Select()
